# Update "想去人数" (F column) counts on the 展览 (sheet 1) and 全部类型 (sheet 4) sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 557
$ws1.Range("F7").Value = 1741
$ws1.Range("F8").Value = 40
$ws1.Range("F10").Value = 138
$ws1.Range("F11").Value = 1879
$ws1.Range("F18").Value = 5
$ws1.Range("F21").Value = 39
$ws1.Range("F23").Value = 1013
$ws1.Range("F24").Value = 3
$ws1.Range("F27").Value = 255
$ws1.Range("F28").Value = 273

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 557
$ws4.Range("F7").Value = 1741
$ws4.Range("F9").Value = 40
$ws4.Range("F11").Value = 138
$ws4.Range("F12").Value = 1879
$ws4.Range("F19").Value = 5
$ws4.Range("F22").Value = 39
$ws4.Range("F24").Value = 1013
$ws4.Range("F25").Value = 3
$ws4.Range("F28").Value = 255
$ws4.Range("F29").Value = 273
